$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the "Updated symbol list" refresh.
# Values are written as text (matching the workbook's existing inlineStr/text
# storage for these columns) by forcing a Text number format before the
# assignment and then clearing the format back to the sheet default so no
# stray style index is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '248.51'
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.58'
$ws.Range("D3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.389'
$ws.Range("D4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05713'
$ws.Range("D5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.409'
$ws.Range("D6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.324'
$ws.Range("D7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8122'
$ws.Range("D8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9258'
$ws.Range("D9").ClearFormats()

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("B10").ClearFormats()

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("C10").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1423'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("E10").ClearFormats()

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("B11").ClearFormats()

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("C11").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07461'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("E11").ClearFormats()

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("B12").ClearFormats()

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("C12").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03125'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("E12").ClearFormats()

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("B13").ClearFormats()

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("C13").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03027'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("E13").ClearFormats()

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("B14").ClearFormats()

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("C14").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09349'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("E14").ClearFormats()

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'MCDex'
$ws.Range("B15").ClearFormats()

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("C15").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.725'
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("E15").ClearFormats()

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("B16").ClearFormats()

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("C16").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001573'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("E16").ClearFormats()

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("B17").ClearFormats()

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("C17").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04768'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("E17").ClearFormats()

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'UpBots'
$ws.Range("B18").ClearFormats()

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("C18").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.01829'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '17UpBotsUBXTBestin24h'
$ws.Range("E18").ClearFormats()

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'One'
$ws.Range("B19").ClearFormats()

$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("C19").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0005792'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '18OneONE'
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.006449'
$ws.Range("D20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.005009'
$ws.Range("D21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.001025'
$ws.Range("D22").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.163'
$ws.Range("D25").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1306'
$ws.Range("D27").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03991'
$ws.Range("D40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006894'
$ws.Range("D41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1068'
$ws.Range("D42").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007538'
$ws.Range("D44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005892'
$ws.Range("D45").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5001'
$ws.Range("D47").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01010'
$ws.Range("D50").ClearFormats()
